# "Generate Report for Handback" — fills in the Latest Target File / Latest
# Handback File / Latest Handback DateTime columns on the per-locale sheets
# after a handback xliff has been processed, and flips the Overview status
# from "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

# Widened column constant (chars). The exact stored XML width of
# 29.9777047293527 isn't reachable through the ColumnWidth char quantization
# (steps of 1/6 char), so use the nearest attainable value (stored width 30).
$wideColWidth = 29.166666666666664
# Column width that serializes to exactly 40.
$fortyColWidth = 39.16666666666667

# ---------------------------------------------------------------------
# Overview sheet: status flips from "Ready for handoff" to the in-sync
# handback message for both rows, and the two status columns (E, F) widen.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsOverview.Columns.Item(5).ColumnWidth = $wideColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideColWidth

# ---------------------------------------------------------------------
# zh-cn sheet: row2 = 4bfa4aa4..., row3 = c775bb97...
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# "Status" column shares the same underlying string as the Overview sheet's
# handoff-status cells, so it flips to the same handed-back message too.
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh.Columns.Item(3).ColumnWidth = $wideColWidth
$wsZh.Columns.Item(9).ColumnWidth = $fortyColWidth
$wsZh.Columns.Item(10).ColumnWidth = $fortyColWidth

$wsZh.Range("J2").Value = "4bfa4aa4-b631-4da9-809c-9931e0a1b6d4.c5aa6656387b8e17046915f57ad0329c7ba78b36.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-27 08:47:36"

$wsZh.Range("J3").Value = "c775bb97-3765-4f9e-8e6f-8678d9e40cd5.f01bfc254bf4b80938d61de26063432e6c36d838.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-27 08:47:36"

# Rebuild the hyperlink set so the new "Latest Target File" links land in
# the same relative order Excel would emit them (A2, I2, A3, I3).
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/74e593e5847fce130625199471bbcf4fed7b6f7f/e2e/4bfa4aa4-b631-4da9-809c-9931e0a1b6d4.md", "", "", "4bfa4aa4-b631-4da9-809c-9931e0a1b6d4.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/74e593e5847fce130625199471bbcf4fed7b6f7f/e2e/4bfa4aa4-b631-4da9-809c-9931e0a1b6d4.md", "", "", "4bfa4aa4-b631-4da9-809c-9931e0a1b6d4.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/74e593e5847fce130625199471bbcf4fed7b6f7f/e2e/c775bb97-3765-4f9e-8e6f-8678d9e40cd5.md", "", "", "c775bb97-3765-4f9e-8e6f-8678d9e40cd5.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/74e593e5847fce130625199471bbcf4fed7b6f7f/e2e/c775bb97-3765-4f9e-8e6f-8678d9e40cd5.md", "", "", "c775bb97-3765-4f9e-8e6f-8678d9e40cd5.md") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet: row2 = 4bfa4aa4..., row3 = c775bb97...
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe.Columns.Item(3).ColumnWidth = $wideColWidth
$wsDe.Columns.Item(9).ColumnWidth = $fortyColWidth
$wsDe.Columns.Item(10).ColumnWidth = $fortyColWidth

$wsDe.Range("J2").Value = "4bfa4aa4-b631-4da9-809c-9931e0a1b6d4.c5aa6656387b8e17046915f57ad0329c7ba78b36.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-27 08:47:42"

$wsDe.Range("J3").Value = "c775bb97-3765-4f9e-8e6f-8678d9e40cd5.f01bfc254bf4b80938d61de26063432e6c36d838.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-27 08:47:42"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/74e593e5847fce130625199471bbcf4fed7b6f7f/e2e/4bfa4aa4-b631-4da9-809c-9931e0a1b6d4.md", "", "", "4bfa4aa4-b631-4da9-809c-9931e0a1b6d4.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/74e593e5847fce130625199471bbcf4fed7b6f7f/e2e/4bfa4aa4-b631-4da9-809c-9931e0a1b6d4.md", "", "", "4bfa4aa4-b631-4da9-809c-9931e0a1b6d4.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/74e593e5847fce130625199471bbcf4fed7b6f7f/e2e/c775bb97-3765-4f9e-8e6f-8678d9e40cd5.md", "", "", "c775bb97-3765-4f9e-8e6f-8678d9e40cd5.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/74e593e5847fce130625199471bbcf4fed7b6f7f/e2e/c775bb97-3765-4f9e-8e6f-8678d9e40cd5.md", "", "", "c775bb97-3765-4f9e-8e6f-8678d9e40cd5.md") | Out-Null
